$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ B=1.02; C=1.032529617257471; D=1.037528786492637; E=0.992614727750844; F=1.042986474742323; I=1.039144001547973; J=1.037659334370648; K=1.040319412124345; L=0.9955398523335997; M=1.04576161424453; N=1.039132929646644 }
    3 = @{ B=1.02; C=1.03339525318935; D=1.038197799438952; E=0.9936372048519299; F=1.044118408323895; I=1.03939614841432; J=1.038167729428411; K=1.040798624205863; L=0.9963617723202687; M=1.046703651088362; N=1.039642046683702 }
    4 = @{ B=1.02; C=1.033955458940294; D=1.03863066992763; E=0.9942998659930998; F=1.044851300271604; I=1.039557984730848; J=1.03849613209115; K=1.041108000084177; L=0.9968940712668347; M=1.04731308511038; N=1.03997091571589 }
    5 = @{ B=1.02; C=1.034190988009161; D=1.038812641332717; E=0.994578699834602; F=1.045159516813173; I=1.039625704315599; J=1.038634057055776; K=1.041237891776644; L=0.9971179600053012; M=1.04756926050467; N=1.040109036549786 }
    6 = @{ B=1.02; C=1.034230535438469; D=1.038843194666068; E=0.994625531979634; F=1.045211274098653; I=1.039637056163441; J=1.038657207295847; K=1.041259691174978; L=0.9971555583673455; M=1.047612271673215; N=1.040132219665854 }
    7 = @{ B=1.02; C=1.033958606019419; D=1.038633101468856; E=0.994303590798249; F=1.044855418249003; I=1.03955889084712; J=1.038497975584697; K=1.041109736372054; L=0.9968970624462089; M=1.047316508259428; N=1.03997276182741 }
    8 = @{ B=1.02; C=1.03282214556159; D=1.037754887126851; E=0.9929600610674297; F=1.043368922606897; I=1.039229488689026; J=1.037831264971538; K=1.040481509861331; L=0.9958175282591056; M=1.046080006075984; N=1.039305104408708 }
    9 = @{ B=1.02; C=1.03082022103726; D=1.036207218369327; E=0.9906006454969559; F=1.040753023025616; I=1.038638955462528; J=1.036652157056976; K=1.039369120322209; L=0.9939188001724441; M=1.043900173391334; N=1.038124322025661 }
    10 = @{ B=1.02; C=1.029486104016188; D=1.03517541409344; E=0.989033133672735; F=1.039011449647385; I=1.038238514950328; J=1.035863251531616; K=1.038623961405404; L=0.9926553831429383; M=1.042446323059584; N=1.037334296163966 }
    11 = @{ B=1.02; C=1.028908545680963; D=1.034728640150103; E=0.988355674866747; F=1.03825788985083; I=1.038063524138836; J=1.035520983387308; K=1.038300463255754; L=0.9921088820399291; M=1.041816643165579; N=1.036991541959634 }
    12 = @{ B=1.02; C=1.028694034091757; D=1.034562690178108; E=0.9881042295826724; F=1.037978067234994; I=1.037998285182272; J=1.03539375022954; K=1.038180176258566; L=0.9919059725120875; M=1.041582729105644; N=1.036864128116189 }
    13 = @{ B=1.02; C=1.028740046670083; D=1.034598286912308; E=0.9881581567098651; F=1.038038086378508; I=1.038012289994602; J=1.035421046667174; K=1.038205983889776; L=0.9919494934313052; M=1.041632905500806; N=1.036891463317896 }
    14 = @{ B=1.02; C=1.028890813687942; D=1.034714922631076; E=0.9883348863814464; F=1.0382347579322; I=1.038058136351231; J=1.035510468281374; K=1.038290522848538; L=0.9920921077337197; M=1.041797308223248; N=1.036981011921044 }
    15 = @{ B=1.02; C=1.028983708803703; D=1.034786786014682; E=0.9884438009545853; F=1.038355944770748; I=1.0380863520708; J=1.035565550720168; K=1.038342593498911; L=0.9921799884222134; M=1.041898599133406; N=1.037036172583219 }
    16 = @{ B=1.02; C=1.029524437322748; D=1.035205065193181; E=0.9890781214508737; F=1.039061472576869; I=1.038250094893858; J=1.035885952747515; K=1.038645413289138; L=0.9926916645766087; M=1.042488109622344; N=1.037357029618196 }
    17 = @{ B=1.02; C=1.029863655563767; D=1.03546744258263; E=0.989476357848556; F=1.039504179551657; I=1.038352379058919; J=1.036086754404718; K=1.038835139924204; L=0.9930127773699352; M=1.042857853050635; N=1.037558116436784 }
    18 = @{ B=1.02; C=1.030061527652627; D=1.03562048319552; E=0.9897087662937556; F=1.03976245639555; I=1.038411885525309; J=1.036203814229161; K=1.038945723219482; L=0.9932001317071769; M=1.043073503357978; N=1.037675342499606 }
    19 = @{ B=1.02; C=1.030128998903588; D=1.035672666109102; E=0.9897880325774034; F=1.039850531175893; I=1.038432149532229; J=1.036243717663407; K=1.038983415497151; L=0.9932640239640975; M=1.043147032030378; N=1.037715302601305 }
    20 = @{ B=1.02; C=1.029827259427631; D=1.03543929193038; E=0.9894336180360679; F=1.039456675764302; I=1.038341420878109; J=1.036065216942119; K=1.038814792428939; L=0.9929783193494215; M=1.042818184607871; N=1.037536548388518 }
    21 = @{ B=1.02; C=1.028846416020027; D=1.034680576261121; E=0.9882828385668249; F=1.038176840769348; I=1.038044642352055; J=1.035484138587; K=1.038265631687502; L=0.9920501090198102; M=1.04174889635082; N=1.036954644835484 }
    22 = @{ B=1.02; C=1.0282298328518; D=1.034203552371548; E=0.9875604150241495; F=1.037372638169754; I=1.037856659997904; J=1.035118216039592; K=1.037919627481109; L=0.9914670000341481; M=1.041076459876144; N=1.036588202636087 }
    23 = @{ B=1.02; C=1.028556684373409; D=1.034456430359411; E=0.9879432794643023; F=1.03779891565199; I=1.037956444229837; J=1.035312252839411; K=1.038103119339806; L=0.991776070289318; M=1.041432943703242; N=1.036782514990418 }
    24 = @{ B=1.02; C=1.029843705238377; D=1.035452011997511; E=0.9894529299347244; F=1.039478140516018; I=1.038346372884657; J=1.036074948993367; K=1.038823986836649; L=0.9929938892766442; M=1.042836109115191; N=1.037546294260396 }
    25 = @{ B=1.02; C=1.031337682373906; D=1.036607337445389; E=0.9912096547607049; F=1.041428880556332; I=1.038792814580899; J=1.036957487305952; K=1.039657331656603; L=0.9944092447426414; M=1.04446382431318; N=1.038430085878613 }
}

foreach ($row in $data.Keys) {
    foreach ($col in $data[$row].Keys) {
        $ws.Range("$col$row").Value = $data[$row][$col]
    }
}
